$d = $word.ActiveDocument

# --- Revert "Justifying Why Stopwords had to be..." ---
#
# The document currently has, in order:
#   [Heading] Not eliminating stop words when indexing
#   [Para]    The elimination of stop-words had to be removed ...
#   [Para]    (contains only the _GoBack bookmark)
#   [Heading] Separating the query processor from the web UI
#   [Para]    Honestly, this had to be done ... (w/ lastRenderedPageBreak)
#   [Para]    (empty)
#   [Para]    Even though this solution ...
#   [Para]    (empty)
#   [Heading] Using a VIEW, and not sending the result set as an array over the port
#
# The target (reverted) order drops the first two "stop words" paragraphs
# and the stray bookmark-only paragraph, reattaches the _GoBack bookmark to
# the end of the "Even though this solution ..." paragraph, merges the
# split "Honestly, this had to be done ..." run around the page break away,
# and moves the lastRenderedPageBreak marker onto the "Using a VIEW ..." run.

# 1) Remove the "Not eliminating stop words when indexing" heading paragraph.
$d.Paragraphs.Item(44).Range.Delete() | Out-Null

# 2) Remove the "The elimination of stop-words ..." paragraph (now shifted to 44).
$d.Paragraphs.Item(44).Range.Delete() | Out-Null

# 3) Remove the now-orphaned bookmark-only paragraph (now shifted to 44); the
#    bookmark itself is recreated in step 5 below.
$d.Paragraphs.Item(44).Range.Delete() | Out-Null

# After the three deletions the paragraphs read:
#   44 Separating the query processor from the web UI
#   45 Honestly, this had to be done ...
#   46 (empty)
#   47 Even though this solution ...
#   48 (empty)
#   49 Using a VIEW, and not sending the result set as an array over the port

# 4) Rebuild paragraph 45 ("Honestly, this had to be done ..."): merge the
#    two runs that used to straddle the page break into one run, and drop
#    the <w:lastRenderedPageBreak/> marker entirely.
$xmlHonestly = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Honestly, this had to be done because the query processor is in Java, while the web UI is in HTML and CSS surely, but also in PHP. To solve this problem, a port was setup to allow for communication between the web UI (PHP) and the query processor </w:t></w:r><w:r><w:t>(Java), treating web UI (PHP) as the client and the query processor (Java) as the server.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(45).Range.InsertXML($xmlHonestly) | Out-Null

# 5) Rebuild paragraph 47 ("Even though this solution ..."): keep its three
#    runs and reattach the _GoBack bookmark at the end of the paragraph.
$xmlEvenThough = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Even though this solution was meant to solve the original problem of d</w:t></w:r><w:r><w:t xml:space="preserve">ifferent languages, it made way </w:t></w:r><w:r><w:t>the search engine more memory efficient, though unfortunately, not more time efficient.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(47).Range.InsertXML($xmlEvenThough) | Out-Null

# 6) Rebuild paragraph 49 ("Using a VIEW ..."): keep its heading formatting
#    and text, and add the <w:lastRenderedPageBreak/> marker that used to sit
#    in the "Honestly, this had to be done ..." paragraph.
$xmlUsingView = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Using a VIEW, and not sending the result set as an array over the port</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(49).Range.InsertXML($xmlUsingView) | Out-Null
